$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 12:33"

# --- Updated COVID stats for existing countries (no reordering) ---
# Iran (row 13)
$ws.Range("B13").Value = 101650
$ws.Range("C13").Value = 1680
$ws.Range("D13").Value = 81587
$ws.Range("E13").Value = 13645
$ws.Range("F13").Value = 2735
$ws.Range("G13").Value = 78
$ws.Range("H13").Value = 6418

# Catar (row 31)
$ws.Range("B31").Value = 17972
$ws.Range("C31").Value = 830
$ws.Range("D31").Value = 2070
$ws.Range("E31").Value = 15890

# Rumania (row 37)
$ws.Range("B37").Value = 14107
$ws.Range("C37").Value = 270
$ws.Range("D37").Value = 5788
$ws.Range("E37").Value = 7461

# Australia (row 52)
$ws.Range("D52").Value = 5984
$ws.Range("E52").Value = 794

# Marruecos (row 56)
$ws.Range("B56").Value = 5382
$ws.Range("C56").Value = 163
$ws.Range("D56").Value = 1969
$ws.Range("E56").Value = 3231
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 182

# Uzbekistan (row 72)
$ws.Range("D72").Value = 1556
$ws.Range("E72").Value = 651

# Eslovenia (row 86)
$ws.Range("B86").Value = 1448
$ws.Range("C86").Value = 3
$ws.Range("D86").Value = 246
$ws.Range("F86").Value = 14
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 99

# Albania (row 98)
$ws.Range("B98").Value = 832
$ws.Range("C98").Value = 12
$ws.Range("D98").Value = 595
$ws.Range("E98").Value = 206

# --- Reorder "Etiopia" ahead of "Birmania" in the countries list, with
#     fresh data for Etiopia, and shift Birmania/Guadalupe/Madagascar
#     down one row each (their own data rows move with them) ---
$ws.Rows(143).Insert()
$ws.Range("A143").Value = "Etiopia"
$ws.Range("B143").Value = 162
$ws.Range("C143").Value = 17
$ws.Range("D143").Value = 93
$ws.Range("E143").Value = 65
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 4

# Remove the old "Etiopia" row (now pushed down to row 147 by the insert
# above) so the total row count — and Gibraltar's row position — is
# unchanged.
$ws.Rows(147).Delete()
